# Fix payload mapping for rbi_email and tiu_email
# Adds explanatory labels in front of the two employee email cells (G4/G5)
# so it is clear which email is the RBI email and which is the TIU email.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F4: label in front of the RBI email (G4 = pgudipati1@rbi.com)
$ws.Range("F4").Value = "                          RBI Email :"
$ws.Range("F4").HorizontalAlignment = -4108   # xlHAlignCenter
$ws.Range("F4").IndentLevel = 0

# F5: label in front of the TIU email (G5 = gudipati.babu@tiuconsulting.com)
$ws.Range("F5").Value = "                                   TIU Email :"
$ws.Range("F5").HorizontalAlignment = -4108   # xlHAlignCenter
$ws.Range("F5").IndentLevel = 0
$ws.Range("F5").VerticalAlignment = -4107     # xlVAlignBottom (clear inherited "top")

# Move the active selection to J7 (matches the saved view state in the workbook)
$ws.Range("J7").Select() | Out-Null
